$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Update the "Model/Script Status" column (G) for several modules
$ws.Range("G21").Value = "y"
$ws.Range("G23").Value = "y"
$ws.Range("G24").Value = ""
$ws.Range("G26").Value = "Script Done. Need Model"
$ws.Range("G30").Value = "Script Done. Need Model"
$ws.Range("G32").Value = ""
$ws.Range("G35").Value = "Basic Script Ready"

# Update the active selection shown when the sheet is opened
$ws.Range("G44").Select()
